$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Devin Booker","PG,SG","Phoenix Suns"),
    @("Immanuel Quickley","PG,SG","Toronto Raptors"),
    @("Jalen Brunson","PG","New York Knicks"),
    @("Desmond Bane","SG,SF","Memphis Grizzlies"),
    @("Myles Turner","C","Indiana Pacers"),
    @("Nick Richards","C","Phoenix Suns"),
    @("Christian Braun","SG,SF","Denver Nuggets"),
    @("Nikola Jovic","PF,C","Miami Heat"),
    @("Walker Kessler","C","Utah Jazz"),
    @("Jalen Williams","SG,SF,PF,C","Oklahoma City Thunder"),
    @("Jimmy Butler","SF,PF","Golden State Warriors"),
    @("Max Christie","SG,SF","Dallas Mavericks"),
    @("Trae Young","PG","Atlanta Hawks"),
    @("Norman Powell","SG,SF","LA Clippers"),
    @("Kawhi Leonard","SG,SF,PF","LA Clippers"),
    @("Brandon Ingram","SG,SF,PF","Toronto Raptors"),
    @("LeBron James","SF,PF","Los Angeles Lakers")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
